$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new LeetCode entry for "Symmetric Tree" on row 22 (day 20).
# Create the hyperlink first (with TextToDisplay set to the address) so the
# hyperlink's "display" attribute matches the target URL text, then
# overwrite the cell's visible text with the problem name afterwards -
# this keeps shared-string allocation order matching (name before url).
$ws.Hyperlinks.Add($ws.Range("B22"), "https://leetcode.com/problems/symmetric-tree/", "", "", "https://leetcode.com/problems/symmetric-tree/") | Out-Null
$ws.Range("B22").Value = "Symmetric Tree"

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 45
$ws.Range("F22").Value = 0.72
$ws.Range("G22").Value = 16.19
$ws.Range("H22").Value = 0.011
$ws.Range("I22").Value = "https://leetcode.com/problems/symmetric-tree/submissions/1062657334/"

# Match the selection left by the author after entering the new row.
$null = $ws.Range("E26").Select()
